# Apply updated simulated-game transition-matrix probabilities to Sheet1.
# Each row represents a starting state; columns hold transition fractions
# (observed-count / total simulated games for that row). More games were
# simulated, so both the counts and the row totals changed, shifting the
# resulting probabilities slightly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: total simulated games = 227
$ws.Range("B2").Value = 38 / 227
$ws.Range("C2").Value = 140 / 227
$ws.Range("J2").Value = 4 / 227
$ws.Range("P2").Value = 26 / 227
$ws.Range("S2").Value = 19 / 227

# Row 3: total simulated games = 141
$ws.Range("B3").Value = 1 / 141
$ws.Range("C3").Value = 4 / 141
$ws.Range("J3").Value = 5 / 141
$ws.Range("P3").Value = 113 / 141
$ws.Range("S3").Value = 18 / 141

# Row 4: total simulated games = 47
$ws.Range("J4").Value = 2 / 47
$ws.Range("P4").Value = 35 / 47
$ws.Range("S4").Value = 10 / 47

# Row 6: total simulated games = 224
$ws.Range("B6").Value = 11 / 224
$ws.Range("D6").Value = 3 / 224
$ws.Range("F6").Value = 10 / 224
$ws.Range("J6").Value = 53 / 224
$ws.Range("O6").Value = 8 / 224
$ws.Range("Q6").Value = 41 / 224
$ws.Range("R6").Value = 15 / 224
$ws.Range("S6").Value = 83 / 224

# Row 7: total simulated games = 155
$ws.Range("B7").Value = 17 / 155
$ws.Range("D7").Value = 3 / 155
$ws.Range("F7").Value = 10 / 155
$ws.Range("J7").Value = 18 / 155
$ws.Range("O7").Value = 1 / 155
$ws.Range("Q7").Value = 33 / 155
$ws.Range("R7").Value = 10 / 155
$ws.Range("S7").Value = 63 / 155

# Row 8: total simulated games = 461
$ws.Range("B8").Value = 46 / 461
$ws.Range("D8").Value = 6 / 461
$ws.Range("F8").Value = 33 / 461
$ws.Range("J8").Value = 51 / 461
$ws.Range("O8").Value = 9 / 461
$ws.Range("Q8").Value = 69 / 461
$ws.Range("R8").Value = 52 / 461
$ws.Range("S8").Value = 195 / 461

# Row 9: total simulated games = 209
$ws.Range("B9").Value = 19 / 209
$ws.Range("D9").Value = 4 / 209
$ws.Range("F9").Value = 21 / 209
$ws.Range("J9").Value = 18 / 209
$ws.Range("O9").Value = 6 / 209
$ws.Range("Q9").Value = 34 / 209
$ws.Range("R9").Value = 23 / 209
$ws.Range("S9").Value = 84 / 209

# Row 10: total simulated games = 1199
$ws.Range("B10").Value = 95 / 1199
$ws.Range("D10").Value = 32 / 1199
$ws.Range("F10").Value = 86 / 1199
$ws.Range("J10").Value = 147 / 1199
$ws.Range("O10").Value = 14 / 1199
$ws.Range("Q10").Value = 269 / 1199
$ws.Range("R10").Value = 145 / 1199
$ws.Range("S10").Value = 411 / 1199

# Row 11: total simulated games = 242
$ws.Range("G11").Value = 44 / 242
$ws.Range("J11").Value = 18 / 242
$ws.Range("K11").Value = 61 / 242
$ws.Range("L11").Value = 111 / 242
$ws.Range("S11").Value = 8 / 242

# Row 12: total simulated games = 113
$ws.Range("G12").Value = 89 / 113
$ws.Range("J12").Value = 16 / 113
$ws.Range("K12").Value = 1 / 113
$ws.Range("L12").Value = 3 / 113
$ws.Range("S12").Value = 4 / 113

# Row 13: total simulated games = 37
$ws.Range("G13").Value = 27 / 37
$ws.Range("J13").Value = 7 / 37
$ws.Range("S13").Value = 3 / 37

# Row 15: total simulated games = 197
$ws.Range("F15").Value = 3 / 197
$ws.Range("H15").Value = 34 / 197
$ws.Range("I15").Value = 16 / 197
$ws.Range("J15").Value = 72 / 197
$ws.Range("K15").Value = 11 / 197
$ws.Range("M15").Value = 2 / 197
$ws.Range("O15").Value = 13 / 197
$ws.Range("S15").Value = 46 / 197

# Row 16: total simulated games = 171
$ws.Range("F16").Value = 1 / 171
$ws.Range("H16").Value = 28 / 171
$ws.Range("I16").Value = 13 / 171
$ws.Range("J16").Value = 77 / 171
$ws.Range("K16").Value = 12 / 171
$ws.Range("M16").Value = 10 / 171
$ws.Range("O16").Value = 8 / 171
$ws.Range("S16").Value = 22 / 171

# Row 17: total simulated games = 443
$ws.Range("F17").Value = 10 / 443
$ws.Range("H17").Value = 72 / 443
$ws.Range("I17").Value = 47 / 443
$ws.Range("J17").Value = 191 / 443
$ws.Range("K17").Value = 30 / 443
$ws.Range("M17").Value = 4 / 443
$ws.Range("O17").Value = 26 / 443
$ws.Range("S17").Value = 63 / 443

# Row 18: total simulated games = 244
$ws.Range("F18").Value = 3 / 244
$ws.Range("H18").Value = 54 / 244
$ws.Range("I18").Value = 21 / 244
$ws.Range("J18").Value = 105 / 244
$ws.Range("K18").Value = 19 / 244
$ws.Range("M18").Value = 3 / 244
$ws.Range("O18").Value = 12 / 244
$ws.Range("S18").Value = 27 / 244

# Row 19: total simulated games = 1193
$ws.Range("F19").Value = 18 / 1193
$ws.Range("H19").Value = 277 / 1193
$ws.Range("I19").Value = 112 / 1193
$ws.Range("J19").Value = 432 / 1193
$ws.Range("K19").Value = 108 / 1193
$ws.Range("M19").Value = 20 / 1193
$ws.Range("N19").Value = 1 / 1193
$ws.Range("O19").Value = 77 / 1193
$ws.Range("S19").Value = 148 / 1193
